# Add a Totals Row to the TEST3 table (xl/tables/table3.xml / sheet3) so that:
#  - the table grows from A1:F6 to A1:F7 and gets a totals row
#  - column "License plate" shows the totals-row label "Összeg"
#  - column "Manufacturing date" / "Price" show MAX, "Is ready for traffic?" shows COUNT
#
# Implemented by writing the totals row content directly (label text in A7,
# SUBTOTAL formulas in C7/D7/F7) and flipping ShowTotals on the table's
# ListObject, which is how Excel itself represents/produces a "Total Row".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST3")

$lo = $ws.ListObjects.Item(1)

# Reuse the date style already used by the "Manufacturing date" column (C2:C5)
# for the new totals-row cell C7, instead of creating a brand new style.
$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Turn on the table's totals row (this is what extends ref to A1:F7 and adds
# totalsRowCount="1" while keeping the autoFilter range at A1:F6).
$lo.ShowTotals = $true

# Totals-row label for the first column.
$ws.Range("A7").Value = "Összeg"

# Totals-row aggregate formulas (MAX / MAX / COUNT), referencing the table's
# structured column names, matching what Excel generates for a Total Row.
$ws.Range("C7").Formula = "=SUBTOTAL(104,TEST3[Manufacturing date])"
$ws.Range("D7").Formula = "=SUBTOTAL(104,TEST3[Price])"
$ws.Range("F7").Formula = "=SUBTOTAL(103,TEST3[Is ready for traffic?])"

# Match the author's final selection in the sheet.
[void]$ws.Range("E12").Select()

# Best-effort: also record the semantic totals-row metadata through the
# ListObject/ListColumns API (label + aggregate function per column), in
# addition to the literal cell content written above.
$lo.ListColumns.Item(1).TotalsRowLabel = "Összeg"
$lo.ListColumns.Item(3).TotalsRowFunction = "max"
$lo.ListColumns.Item(4).TotalsRowFunction = "max"
$lo.ListColumns.Item(6).TotalsRowFunction = "count"
